# ShowtimeExport.xlsx template update:
#   - rename sheet "DANH SACH PHIM" -> "Sheet1"
#   - retitle the report header
#   - insert a new "Gia ve" column ahead of the existing "Bi xoa" column
#   - relabel the ID-based headers as friendlier, name-based headers
#   - resize the data columns to the new template's layout
#   - move the active selection the way the saved template has it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Sheet tab name.
$ws.Name = "Sheet1"

# 2. Make room for the new "Gia ve" (ticket price) column: inserting at the
#    old "Bi xoa" column (E) shifts it right to F and gives us a blank E to
#    fill in below. Column F therefore keeps its original width/formatting.
$ws.Columns("E").Insert()

# 3. Report title (row 1, merged A1:F1).
$ws.Range("A1").Value = "DANH SÁCH SUẤT CHIẾU"

# 4. Header row (row 3).
$ws.Range("C3").Value = "Rạp chiếu"
$ws.Range("B3").Value = "Tên Phim"
# D3 "Giờ chiếu" is unchanged.
$ws.Range("E3").Value = "Giá vé"
# F3 "Bị xóa" shifted in place automatically by the column insert above.

# 5. Column widths for the new layout. Excel quantises ColumnWidth to whole
#    pixels against the workbook's Normal-style Maximum Digit Width (7px for
#    this template's Times New Roman 11 default font), so the inputs below
#    are chosen to land on the pixel closest to the template's target
#    character widths (40.77734375 / 32.44140625 / 40.44140625 / 17.88671875).
$ws.Columns("B").ColumnWidth = 40
$ws.Columns("C").ColumnWidth = 31.714285714285715
$ws.Columns("D").ColumnWidth = 39.714285714285715
$ws.Columns("E").ColumnWidth = 17.142857142857146
# Column F keeps the original "Bi xoa" column width (17.21875) untouched.

# 6. Active cell/selection as saved in the template.
[void]$ws.Range("F13").Select()

Write-Output "done"
